$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- Shared-string text edits (Volume/Number + report date range) ---
$c1 = $ws.Range("A8").Characters(21, 2)
$c1.Text = "16"
$c3 = $ws.Range("C9").Characters(46, 9)
$c3.Text = "4/20/2025"
$c2 = $ws.Range("C9").Characters(27, 8)
$c2.Text = "4/14/2025"

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 12
$ws.Range("K15").Value = 300
$ws.Range("L15").Value = 500
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 0
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 64
$ws.Range("J16").Value = 61
$ws.Range("K16").Value = 4.918032786885
$ws.Range("L16").Value = 68.421052631578
$ws.Range("M16").Value = 77.777777777777
$ws.Range("N16").Value = -83.068783068783
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -69.230769230769
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = -19.354838709677
$ws.Range("I17").Value = 93
$ws.Range("J17").Value = 73
$ws.Range("K17").Value = 27.397260273972
$ws.Range("L17").Value = 52.459016393442
$ws.Range("M17").Value = 151.351351351351
$ws.Range("N17").Value = -36.301369863013
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 7.692307692307
$ws.Range("I18").Value = 56
$ws.Range("J18").Value = 56
$ws.Range("L18").Value = 33.333333333333
$ws.Range("M18").Value = -5.084745762711
$ws.Range("N18").Value = -91.397849462365
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 38
$ws.Range("E19").Value = -42.105263157894
$ws.Range("F19").Value = 122
$ws.Range("G19").Value = 134
$ws.Range("H19").Value = -8.955223880597
$ws.Range("I19").Value = 540
$ws.Range("J19").Value = 537
$ws.Range("K19").Value = 0.558659217877
$ws.Range("L19").Value = -1.818181818181
$ws.Range("M19").Value = 17.903930131004
$ws.Range("N19").Value = -74.100719424460
$ws.Range("I20").Value = 19
$ws.Range("J20").Value = 12
$ws.Range("K20").Value = 58.333333333333
$ws.Range("L20").Value = -9.523809523809
$ws.Range("M20").Value = 72.727272727272
$ws.Range("N20").Value = -87.581699346405
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 60
$ws.Range("E21").Value = -38.333333333333
$ws.Range("F21").Value = 178
$ws.Range("G21").Value = 194
$ws.Range("H21").Value = -8.247422680412
$ws.Range("I21").Value = 784
$ws.Range("J21").Value = 743
$ws.Range("K21").Value = 5.518169582772
$ws.Range("L21").Value = 9.803921568627
$ws.Range("M21").Value = 28.104575163398
$ws.Range("N21").Value = -77.142857142857
$ws.Range("D22").Value = 2
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = -80
$ws.Range("I22").Value = 21
$ws.Range("J22").Value = 28
$ws.Range("K22").Value = -25
$ws.Range("L22").Value = -4.545454545454
$ws.Range("M22").Value = 10.526315789473
$ws.Range("C24").Value = 55
$ws.Range("D24").Value = 50
$ws.Range("E24").Value = 10
$ws.Range("F24").Value = 217
$ws.Range("G24").Value = 182
$ws.Range("H24").Value = 19.230769230769
$ws.Range("I24").Value = 785
$ws.Range("J24").Value = 839
$ws.Range("K24").Value = -6.436233611442
$ws.Range("L24").Value = 6.802721088435
$ws.Range("M24").Value = 53.620352250489
$ws.Range("C25").Value = 45
$ws.Range("D25").Value = 56
$ws.Range("E25").Value = -19.642857142857
$ws.Range("F25").Value = 187
$ws.Range("G25").Value = 184
$ws.Range("H25").Value = 1.630434782608
$ws.Range("I25").Value = 753
$ws.Range("J25").Value = 836
$ws.Range("K25").Value = -9.928229665071
$ws.Range("L25").Value = -1.697127937336
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 24
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 49
$ws.Range("G26").Value = 60
$ws.Range("H26").Value = -18.333333333333
$ws.Range("I26").Value = 194
$ws.Range("J26").Value = 197
$ws.Range("K26").Value = -1.522842639593
$ws.Range("L26").Value = 1.570680628272
$ws.Range("M26").Value = 32.876712328767
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("I27").Value = 13
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = 116.666666666667
$ws.Range("L27").Value = 116.666666666667
$ws.Range("C28").Value = 3
$ws.Range("E28").Value = -25
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 38
$ws.Range("J28").Value = 22
$ws.Range("K28").Value = 72.727272727272
$ws.Range("L28").Value = 72.727272727272

# --- Shared-string -> numeric conversions (also fix number format/style) ---
$ws.Range("C15").Value = 2
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("C20").Value = 1
$ws.Range("C20").NumberFormat = '#,##0'
$ws.Range("D20").Value = 1
$ws.Range("D20").NumberFormat = '#,##0'
$ws.Range("E20").Value = 0
$ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C27").Value = 2
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = 100
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'

# --- Numeric -> shared-string conversions (force text, then restore style via format paste) ---
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "***.*"
$ws.Range("E23").Copy()
$ws.Range("E31").PasteSpecial(-4122)
